$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh (GitHub Actions bot): update the Price (D) and
# Volume(1h) (E) columns with the latest scraped figures.
#
# Price values that look like a plain decimal number (e.g. "0.9978" or
# "6.100") are written with a leading apostrophe - the normal Excel way
# to force text entry - so they stay exact text (keeping trailing zeros)
# instead of silently being re-interpreted as numbers. Values that are
# already not valid numbers (e.g. "26.792.21") do not need this.

$ws.Range('D2').Value = '26.792.21'
$ws.Range('E2').Value = '  +4.14%  '
$ws.Range('D3').Value = '1.865.89'
$ws.Range('E3').Value = '  +2.61%  '
$ws.Range('D4').Value = '''0.9978'
$ws.Range('E4').Value = '  -0.36%  '
$ws.Range('D5').Value = '''273.77'
$ws.Range('E5').Value = '  -1.85%  '
$ws.Range('D6').Value = '''0.9981'
$ws.Range('E6').Value = '  -0.28%  '
$ws.Range('D7').Value = '''0.5272'
$ws.Range('E7').Value = '  +3.66%  '
$ws.Range('D8').Value = '''0.3384'
$ws.Range('E8').Value = '  -4.41%  '
$ws.Range('D9').Value = '''0.06813'
$ws.Range('E9').Value = '  +1.90%  '
$ws.Range('E10').Value = '  -0.26%  '
$ws.Range('D11').Value = '''0.7936'
$ws.Range('E12').Value = '  -1.62%  '
$ws.Range('D13').Value = '1.875.22'
$ws.Range('E13').Value = '  +3.09%  '
$ws.Range('D14').Value = '''89.94'
$ws.Range('E14').Value = '  +2.51%  '
$ws.Range('D15').Value = '''5.131'
$ws.Range('E15').Value = '  +1.00%  '
$ws.Range('D16').Value = '''0.9988'
$ws.Range('E16').Value = '  -0.21%  '
$ws.Range('E17').Value = '  +2.43%  '
$ws.Range('D18').Value = '''0.000008007'
$ws.Range('E18').Value = '  -0.50%  '
$ws.Range('D19').Value = '''0.9987'
$ws.Range('E19').Value = '  -0.22%  '
$ws.Range('D20').Value = '26.768.47'
$ws.Range('E20').Value = '  +3.85%  '
$ws.Range('D21').Value = '2.096.53'
$ws.Range('E21').Value = '  +2.36%  '
$ws.Range('D22').Value = '''4.717'
$ws.Range('E22').Value = '  -0.79%  '
$ws.Range('D23').Value = '''9.982'
$ws.Range('E23').Value = '  -0.23%  '
$ws.Range('D24').Value = '''6.100'
$ws.Range('E24').Value = '  -0.15%  '
$ws.Range('D25').Value = '''2.361'
$ws.Range('E25').Value = '  +5.25%  '
$ws.Range('D26').Value = '''145.61'
$ws.Range('E26').Value = '  +2.29%  '
$ws.Range('D27').Value = '''1.651'
$ws.Range('E27').Value = '  -1.09%  '
$ws.Range('D28').Value = '''17.21'
$ws.Range('D29').Value = '''112.43'
$ws.Range('E29').Value = '  +2.92%  '
$ws.Range('D30').Value = '''4.327'
$ws.Range('E30').Value = '  -0.26%  '
$ws.Range('D31').Value = '''4.311'
$ws.Range('E31').Value = '  +1.99%  '
$ws.Range('D32').Value = '''0.08878'
$ws.Range('E32').Value = '  +1.36%  '
$ws.Range('D33').Value = '''0.04921'
$ws.Range('E33').Value = '  +0.63%  '
$ws.Range('D34').Value = '''1.164'
$ws.Range('E34').Value = '  +2.49%  '
$ws.Range('D35').Value = '''0.7275'
$ws.Range('E35').Value = '  -0.12%  '
$ws.Range('D36').Value = '''2.877'
$ws.Range('E36').Value = '  -0.39%  '
$ws.Range('E37').Value = '  +2.47%  '
$ws.Range('D38').Value = '''2.336'
$ws.Range('E38').Value = '  -0.99%  '
$ws.Range('E39').Value = '  -0.45%  '
$ws.Range('D40').Value = '''0.5095'
$ws.Range('E40').Value = '  -1.34%  '
$ws.Range('D41').Value = '''0.9402'
$ws.Range('E41').Value = '  -3.24%  '
$ws.Range('D42').Value = '''116.06'
$ws.Range('E42').Value = '  +1.57%  '
$ws.Range('D43').Value = '''6.121'
$ws.Range('E43').Value = '  -1.89%  '
$ws.Range('D44').Value = '''7.997'
$ws.Range('E44').Value = '  -0.15%  '
$ws.Range('D45').Value = '''0.9980'
$ws.Range('E45').Value = '  -0.27%  '
$ws.Range('D46').Value = '''0.4414'
$ws.Range('E46').Value = '  -2.72%  '
$ws.Range('E47').Value = '  -3.01%  '
$ws.Range('D48').Value = '''9.276'
$ws.Range('E48').Value = '  +0.82%  '
$ws.Range('D49').Value = '''36.09'
$ws.Range('D50').Value = '''0.05948'
$ws.Range('E50').Value = '  +1.89%  '
$ws.Range('E51').Value = '  -1.75%  '
